$d = $word.ActiveDocument

# Update the date/weekday heading
$d.Content.Find.Execute("2025-10-13 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-10-14 Tuesday", 2)

# Update the practice-problem table (5 data rows x 5 columns living in
# table rows 1,5,9,13,17 of the 20-row table). Addressing by (row,col)
# avoids ambiguity since some old values (e.g. "60÷4=") repeat.
$t = $d.Tables.Item(1)

$newValues = @{
    1  = "23÷4="; 2  = "90÷3="; 3  = "59÷5="; 4  = "64÷2="; 5  = "82÷2=";
    21 = "84÷5="; 22 = "32÷5="; 23 = "64÷7="; 24 = "74÷8="; 25 = "36÷8=";
    41 = "95÷6="; 42 = "40÷5="; 43 = "21÷4="; 44 = "22÷6="; 45 = "25÷3=";
    61 = "40÷9="; 62 = "14÷7="; 63 = "98÷7="; 64 = "16÷5="; 65 = "35÷8=";
    81 = "17÷2="; 82 = "39÷9="; 83 = "58÷4="; 84 = "90÷2="; 85 = "94÷8=";
}

$rows = @(1, 5, 9, 13, 17)
foreach ($r in $rows) {
    for ($c = 1; $c -le 5; $c++) {
        $key = (($r - 1) * 5) + $c
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$key]
    }
}
